$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Fix the autofilter: it was mistakenly filtering Start-time column (H) on
# "12.50"; it should filter the Class column (C) on
# "Fraud and Business Process Analytics" instead. Clearing/reapplying the
# filter also re-evaluates row visibility for every data row (1:1 with the
# hidden/unhidden rows in the target diff, including un-hiding row 76, the
# Fraud and Business Process Analytics class row).
$dataRange = $ws.Range("A1:M99")
$dataRange.AutoFilter(8)
$dataRange.AutoFilter(3, @("Fraud and Business Process Analytics"), 7)

# --- Correct the mistaken value in F76 (Fraud and Business Process Analytics row)
$ws.Range("F76").Value = 0

# --- Update the remembered selection in the frozen (bottom-left) pane
$ws.Range("F100").Select()
